# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The old statement covered 3 workers (EDWIN FERNANDO MACIAS MARTELO,
# ALFREDO ANTONIO NIEVES ATENCIO, MARIA ELENA SOLER CUBEROS) across 9
# periods. The new "parte 1" drops EDWIN entirely, keeps ALFREDO (one
# period) and MARIA (now re-sorted into ascending period order with an
# updated Salario Basico), and refreshes the summary totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove EDWIN FERNANDO MACIAS MARTELO's 7 data rows entirely (rows 16-22).
#    Everything below shifts up automatically (including the closing
#    signature rows and all merged ranges).
$ws.Rows("16:22").Delete()

# After the delete, the data block (rows 16-23) reads:
#   16: ALFREDO  1712
#   17: MARIA    2004
#   18: MARIA    2003
#   19: MARIA    2002
#   20: MARIA    2001
#   21: MARIA    1912
#   22: MARIA    1911
#   23: MARIA    1910
# The refreshed report re-sorts MARIA's periods ascending (1910 -> 2004)
# and raises her Salario Basico from 828200 to 828116.

$periods = @("1910", "1911", "1912", "2001", "2002", "2003", "2004")
$valorMora = @(33128, 33128, 33125, 33125, 33125, 33125, 33125)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 17 + $i
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = $valorMora[$i]
    $ws.Cells.Item($r, 7).Value = 828116
}

# 2) Refresh the summary header figures.
$ws.Range("E11").Value = 243685
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 8
